$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-39 with new Actual Consumption and Timestamp values
$newData = @(
    @(2, 5187, 45859),
    @(3, 5169, 45859.01041666666),
    @(4, 5060, 45859.02083333334),
    @(5, 5033, 45859.03125),
    @(6, 5004, 45859.04166666666),
    @(7, 5004, 45859.05208333334),
    @(8, 5006, 45859.0625),
    @(9, 4925, 45859.07291666666),
    @(10, 4919, 45859.08333333334),
    @(11, 4934, 45859.09375),
    @(12, 4932, 45859.10416666666),
    @(13, 4991, 45859.11458333334),
    @(14, 4985, 45859.125),
    @(15, 4932, 45859.13541666666),
    @(16, 4960, 45859.14583333334),
    @(17, 4986, 45859.15625),
    @(18, 5016, 45859.16666666666),
    @(19, 5019, 45859.17708333334),
    @(20, 5030, 45859.1875),
    @(21, 5042, 45859.19791666666),
    @(22, 5251, 45859.20833333334),
    @(23, 5328, 45859.21875),
    @(24, 5411, 45859.22916666666),
    @(25, 5528, 45859.23958333334),
    @(26, 5732, 45859.25),
    @(27, 5930, 45859.26041666666),
    @(28, 5957, 45859.27083333334),
    @(29, 5995, 45859.28125),
    @(30, 6064, 45859.29166666666),
    @(31, 6087, 45859.30208333334),
    @(32, 6171, 45859.3125),
    @(33, 6093, 45859.32291666666),
    @(34, 6140, 45859.33333333334),
    @(35, 6038, 45859.34375),
    @(36, 6062, 45859.35416666666),
    @(37, 6021, 45859.36458333334),
    @(38, 5959, 45859.375),
    @(39, 5911, 45859.38541666666)
)

foreach ($item in $newData) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}

# Remove now-unused trailing rows 40-47 (previously present, data shrunk)
$ws.Range("A40:B47").EntireRow.Delete()

